{"js": "// Office.js (Word JavaScript API) script.\n// Applies the meaningful content changes described by the diff:\n//   1. \"Usually I have used\"        -> \"Usually, I have used\"          (comma added)\n//   2. \"the video I just followed\"  -> \"the video, I just followed\"    (comma added)\n//   3. Appends a blank spacer paragraph followed by two new bold\n//      Finnish diary paragraphs at the very end of the document body.\n//\n// (The diff also contains a large amount of non-semantic \"save churn\"\n//  produced by Word itself - w:proofErr spell/grammar markers, w:lang\n//  attributes being added/removed, style-sheet lsdException bookkeeping,\n//  and customXml part renumbering. None of that changes the document's\n//  visible text/formatting, so it is intentionally not reproduced here.)\n\nconst body = context.document.body;\n\n// --- 1. \"Usually I have used\" -> \"Usually, I have used\" -----------------\nconst usuallyResults = body.search(\"Usually I have used\", { matchCase: true, matchWholeWord: false });\nusuallyResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < usuallyResults.items.length; i++) {\n  usuallyResults.items[i].insertText(\"Usually, I have used\", Word.InsertLocation.replace);\n}\n\n// --- 2. \"the video I just followed\" -> \"the video, I just followed\" -----\nconst videoResults = body.search(\"the video I just followed\", { matchCase: true, matchWholeWord: false });\nvideoResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < videoResults.items.length; i++) {\n  videoResults.items[i].insertText(\"the video, I just followed\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n\n// --- 3. Append new paragraphs at the end of the document ----------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Blank spacer paragraph (matches the \"Leipteksti\" body-text style used\n// throughout the diary, inherited automatically from the paragraph it is\n// split off from).\nconst spacerParagraph = lastParagraph.insertParagraph(\"\", Word.InsertLocation.after);\n\n// First new diary entry (bold Finnish note).\nconst entryParagraph1 = spacerParagraph.insertParagraph(\n  \"1.12 tein pari parttia eteenp\u00e4in button ja menu overlay\",\n  Word.InsertLocation.after\n);\nentryParagraph1.font.bold = true;\n\n// Second new diary entry (bold Finnish note).\nconst entryParagraph2 = entryParagraph1.insertParagraph(\n  \"ei mit\u00e4\u00e4n erikoista ollu kunha taas tappelin npm kanssa\",\n  Word.InsertLocation.after\n);\nentryParagraph2.font.bold = true;\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the meaningful content changes described by the diff:\n#   1. \"Usually I have used\"        -> \"Usually, I have used\"          (comma added)\n#   2. \"the video I just followed\"  -> \"the video, I just followed\"    (comma added)\n#   3. Appends a blank spacer paragraph followed by two new bold\n#      Finnish diary paragraphs at the very end of the document body.\n#\n# (The diff also contains a large amount of non-semantic \"save churn\"\n#  produced by Word itself - w:proofErr spell/grammar markers, w:lang\n#  attributes being added/removed, style-sheet lsdException bookkeeping,\n#  and customXml part renumbering. None of that changes the document's\n#  visible text/formatting, so it is intentionally not reproduced here.)\n\n$d = $word.ActiveDocument\n\n# --- 1. \"Usually I have used\" -> \"Usually, I have used\" -----------------\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Execute(\"Usually I have used\", $false, $false, $false, $false, $false, $true, 1, $false, \"Usually, I have used\", 2) | Out-Null\n\n# --- 2. \"the video I just followed\" -> \"the video, I just followed\" -----\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\"the video I just followed\", $false, $false, $false, $false, $false, $true, 1, $false, \"the video, I just followed\", 2) | Out-Null\n\n# --- 3. Append new paragraphs at the end of the document ----------------\n\n# Blank spacer paragraph (matches the \"Leipteksti\" body-text style used\n# throughout the diary, inherited automatically from the paragraph mark it\n# is split off from).\n$lastParagraph = $d.Paragraphs.Last\n$endRange = $lastParagraph.Range\n$endRange.Collapse(0) | Out-Null\n$endRange.InsertParagraphAfter()\n\n# First new diary entry (bold Finnish note).\n$spacerParagraph = $d.Paragraphs.Last\n$spacerRange = $spacerParagraph.Range\n$spacerRange.Collapse(0) | Out-Null\n$spacerRange.InsertParagraphAfter()\n\n$entryParagraph1 = $d.Paragraphs.Last\n$entryRange1 = $entryParagraph1.Range\n$entryRange1.Text = \"1.12 tein pari parttia eteenp\u00e4in button ja menu overlay\"\n$entryRange1.Font.Bold = 1\n\n# Second new diary entry (bold Finnish note).\n$entryRange1b = $entryParagraph1.Range\n$entryRange1b.Collapse(0) | Out-Null\n$entryRange1b.InsertParagraphAfter()\n\n$entryParagraph2 = $d.Paragraphs.Last\n$entryRange2 = $entryParagraph2.Range\n$entryRange2.Text = \"ei mit\u00e4\u00e4n erikoista ollu kunha taas tappelin npm kanssa\"\n$entryRange2.Font.Bold = 1\n"}
